$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing changed cells (rows 675-677)
$ws.Range("F675").Value = 3798.03799633
$ws.Range("E676").Value = 29190.55
$ws.Range("F676").Value = 13725.97767151
$ws.Range("F677").Value = 22921.12890833

# Append new rows 678-691 with OHLCV data
$data = @(
    @(45147.41666666666, 29775.65, 30129.27, 29365.49, 29573.89, 15359.03749689),
    @(45148.41666666666, 29573.92, 29712.95, 29317.25, 29433.51, 11115.30974098),
    @(45149.41666666666, 29434.01, 29537.54, 29223.42, 29407.86, 8933.77670339),
    @(45150.41666666666, 29407.86, 29473.73, 29361.72, 29422.34, 2842.16799777),
    @(45151.41666666666, 29422.42, 29451.93, 29264.29, 29289.76, 3248.0312529),
    @(45152.41666666666, 29288.97, 29667.77, 29090.49, 29419.22, 10383.36599966),
    @(45153.41666666666, 29419, 29467.16, 29064.65, 29176.89, 8584.02071459),
    @(45154.41666666666, 29176.63, 29232.71, 28701.67, 28707.5, 14000.5675284),
    @(45155.41666666666, 28705.2, 28758.96, 25253.44, 26635.04, 43569.94033756),
    @(45156.41666666666, 26631.58, 26824.09, 25618.28, 26053.12, 28983.00011015),
    @(45157.41666666666, 26054.35, 26267.78, 25800.8, 26097.91, 8854.32722316),
    @(45158.41666666666, 26096.9, 26295.77, 25987.68, 26196.16, 6240.6993679),
    @(45159.41666666666, 26195.97, 26251.06, 25820.83, 26129.39, 13690.79828458),
    @(45160.41666666666, 26129.39, 26138.6, 25361.73, 26046.38, 16916.00830198)
)

$startRow = 678
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # Copy the datetime cell style (s="2", YYYY-MM-DD HH:MM:SS format) from the row above
    $ws.Range("A677").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
